$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").EntireColumn.Delete()

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "203080756"
$ws.Range("D2").ClearFormats()
